# The underlying change here is a plain "open in Excel, then save" round
# trip: the workbook was originally produced by Apache POI (docProps/app.xml
# -> "Apache POI"), and Excel's own re-serialisation adds its usual
# housekeeping XML (fileVersion/mc:AlternateContent/xr:revisionPtr, a
# calcPr element, the standard Office theme part, x14ac/xr namespace
# declarations + dyDescent/spans row attributes, cellStyles/dxfs/tableStyles
# boilerplate in styles.xml, numeric attribute normalisation such as
# sz val="11.0" -> "11", etc.) without touching a single cell value, shared
# string, formula, or any visible formatting. No row/column/cell in the
# "Details" sheet was added, removed, or changed.
#
# So the faithful reproduction of the commit is simply to open the workbook
# and save it back out, letting the host re-emit the canonical OOXML -
# without making any data or formatting edits of our own.

$wb = $excel.ActiveWorkbook
$wb.Save()
